# Add a new "Spain" market test-data sheet, based on the existing
# "Italy" sheet (same layout/styles), then trim it down to the Spain
# specifics: drop the "MZX Communicator" printer row (not applicable
# for Spain) and fill in the Spain market name + printer/ticket codes.

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Duplicate "Italy" and drop the copy right after it -> becomes the new
# last tab, picks up the next sheetId/rId, and becomes the active tab.
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# Spain-specific market name / ticket references.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2064/T2062/"

# Spain doesn't use the MZX Communicator printer - remove that row,
# shifting everything below it up by one.
$spain.Rows("12:12").Delete() | Out-Null

# Leave the cursor where the edits left off.
$spain.Range("B12").Select() | Out-Null
